$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B5").Value = "SingleUseId1"
$ws.Range("C5").Value = "Default"
$ws.Range("D5").Value = "Left"
$ws.Range("E5").Value = "LTR"
$ws.Range("F5").Value = "99999"
